$wb = $excel.ActiveWorkbook

# --- Sheet "Estadisticos 1P" ---
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("D2").Value = 0
$ws1.Range("E2").Value = 19
$ws1.Range("F2").Value = 14
$ws1.Range("G2").Value = 42.42
$ws1.Range("H2").Value = 6.2

$ws1.Range("D3").Value = 0
$ws1.Range("E3").Value = 8
$ws1.Range("F3").Value = 13
$ws1.Range("G3").Value = 61.9
$ws1.Range("H3").Value = 7.2

$ws1.Range("D4").Value = 0
$ws1.Range("E4").Value = 12
$ws1.Range("F4").Value = 25
$ws1.Range("G4").Value = 67.57
$ws1.Range("H4").Value = 7.7

# --- Sheet "Estadisticos 2P" ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("E2").Value = 33
$ws2.Range("E3").Value = 21
$ws2.Range("E4").Value = 37

# --- Sheet "Estadisticos Final" ---
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("D2").Value = 0
$ws3.Range("E2").Value = 19
$ws3.Range("F2").Value = 14
$ws3.Range("G2").Value = 42.42
$ws3.Range("H2").Value = 6.2

$ws3.Range("D3").Value = 0
$ws3.Range("E3").Value = 8
$ws3.Range("F3").Value = 13
$ws3.Range("G3").Value = 61.9
$ws3.Range("H3").Value = 7.2

$ws3.Range("D4").Value = 0
$ws3.Range("E4").Value = 12
$ws3.Range("F4").Value = 25
$ws3.Range("G4").Value = 67.57
$ws3.Range("H4").Value = 7.7

# --- Sheet "Rescatables" ---
# Remove the two students whose data was deleted from the roster
# (original row 3 = GARCIA ..., original row 5 = RAMIREZ ...)
# Delete bottom-up so row indices stay valid.
$ws4 = $wb.Worksheets.Item("Rescatables")
$ws4.Rows(5).Delete() | Out-Null
$ws4.Rows(3).Delete() | Out-Null
